$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (05-13-2015, Wednesday) gains Official Business Network entries ---
# First, copy the cell formatting from row 10 (which already carries the
# "has OB data" style) onto row 9 so its fill/border/font match (style 11 -> 9)
$ws.Range("A10:P10").Copy() | Out-Null
$ws.Range("A9:P9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Now populate the Official Business Departure/Start/End/Arrival times and remarks
$ws.Range("K9").Value = "08:30:00"
$ws.Range("L9").Value = "08:30:00"
$ws.Range("M9").Value = "18:30:00"
$ws.Range("N9").Value = "18:30:00"
$ws.Range("P9").Value = "~OB Network| R"

# --- Append " R" to the two existing Official Business remarks ---
# Row 10 (05-14-2015, Thursday)
$ws.Range("P10").Value = "~OB Others|Sit| R"

# Rows 14 & 15 (05-18-2015 and 05-19-2015) share the same remark text
$ws.Range("P14").Value = "~OB Others|SIT| R"
$ws.Range("P15").Value = "~OB Others|SIT| R"

# --- Update the manually-entered Total Overtime Hours figure ---
$ws.Range("I23").Value = 3.5
